$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header details ---
$ws.Range("C2").Value = "Hartmut"

# B3 holds a 16-digit card-number string; force text storage (like the
# original inline-string cell) but normalize the style back to the plain
# "s=8" formatting used throughout the sheet (avoid picking up a stray
# text/quote-prefix number format from the coercion).
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null

$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 08.01.2025"

# --- Row 6 (existing transaction, dates/desc/amount updated) ---
$ws.Range("B6").Value = "11.01."
$ws.Range("C6").Value = "12.01."
$ws.Range("D6").Value = "RECHNUNG VODAFONE GMBH 86537976"
$ws.Range("E6").Value = "41,39-"

# --- Row 7 ---
$ws.Range("B7").Value = "13.01."
$ws.Range("C7").Value = "14.01."
$ws.Range("D7").Value = "PAYPAL GOFKKZ"
$ws.Range("E7").Value = "16,94-"

# --- Row 8 ---
$ws.Range("B8").Value = "14.01."
$ws.Range("C8").Value = "15.01."
$ws.Range("D8").Value = "KARTENZ./14.01 REWE RO"
$ws.Range("E8").Value = "63,16-"

# --- Row 9 ---
$ws.Range("B9").Value = "16.01."
$ws.Range("C9").Value = "17.01."
$ws.Range("D9").Value = "KARTENZ./16.01 EDEKA RO"
$ws.Range("E9").Value = "31,64-"

# --- Row 10 (was blank, now a new transaction) ---
# B10/C10/D10 already carry the correct "s=8" format from the template, so
# only the value needs to be written. E10 used the blank-row style (s=12);
# bring it in line with the other amount cells (s=17) by copying formats
# from E9, then restore the value.
$ws.Range("B10").Value = "19.01."
$ws.Range("C10").Value = "20.01."
$ws.Range("D10").Value = "ZALANDO MKTPLC EU IEYHWY"
$ws.Range("E9").Copy() | Out-Null
$ws.Range("E10").PasteSpecial(-4122) | Out-Null
$ws.Range("E10").Value = "162,00-"

# --- Row 11 (was blank, now a new transaction) ---
$ws.Range("B11").Value = "20.01."
$ws.Range("C11").Value = "21.01."
$ws.Range("D11").Value = "KARTENZAHLUNG JET TANKSTELLE"
$ws.Range("E9").Copy() | Out-Null
$ws.Range("E11").PasteSpecial(-4122) | Out-Null
$ws.Range("E11").Value = "77,08-"

# --- Closing balance line ---
$ws.Range("D12").Value = "KONTOSTAND AM 22.01.2025"
$ws.Range("E12").Value = "392,21-"

# --- Next statement date footer ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 29.01.2025"
